$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty I18 cell (becomes a truly blank cell)
$ws.Range("I18").ClearContents()

# Row 19
$ws.Range("A19").Value = '0TF26852'
$ws.Range("B19").Value = 'MONTAGNE MASCARILLA CHUPA CHUPS NARANJA'
$ws.Range("C19").Value = 'TRAT.FEMENINO'
$ws.Range("D19").Value = 'Tiene PT'
$ws.Range("E19").Value = 'Tiene ES'
$ws.Range("F19").Value = 'Tiene IT'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '8'
$ws.Range("H19").Value = 'ML'
$ws.Range("J19").Value = 'Solo Revisión'

# Row 20
$ws.Range("A20").Value = '0TF26855'
$ws.Range("B20").Value = 'MONTAGNE MASCARILLA CHUPA CHUPS MANZANA'
$ws.Range("C20").Value = 'TRAT.FEMENINO'
$ws.Range("D20").Value = 'Tiene PT'
$ws.Range("E20").Value = 'Tiene ES'
$ws.Range("F20").Value = 'Tiene IT'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '8'
$ws.Range("H20").Value = 'ML'
$ws.Range("J20").Value = 'Solo Revisión'

# Row 21
$ws.Range("A21").Value = '0TF26854'
$ws.Range("B21").Value = 'MONTAGNE MASCARILLA CHUPA CHUPS SANDIA'
$ws.Range("C21").Value = 'TRAT.FEMENINO'
$ws.Range("D21").Value = 'Tiene PT'
$ws.Range("E21").Value = 'Tiene ES'
$ws.Range("F21").Value = 'Tiene IT'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '8'
$ws.Range("H21").Value = 'ML'
$ws.Range("J21").Value = 'Solo Revisión'

# Row 22
$ws.Range("A22").Value = '0TF26856'
$ws.Range("B22").Value = 'MONTAGNE MASCARILLA CHUPA CHUPS COLA'
$ws.Range("C22").Value = 'TRAT.FEMENINO'
$ws.Range("D22").Value = 'Tiene PT'
$ws.Range("E22").Value = 'Tiene ES'
$ws.Range("F22").Value = 'Tiene IT'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '8'
$ws.Range("H22").Value = 'ML'
$ws.Range("J22").Value = 'Solo Revisión'

# Row 23
$ws.Range("A23").Value = '0TF26853'
$ws.Range("B23").Value = 'MONTAGNE MASCARILLA CHUPA CHUPS FRESA NATA'
$ws.Range("C23").Value = 'TRAT.FEMENINO'
$ws.Range("D23").Value = 'Tiene PT'
$ws.Range("E23").Value = 'Tiene ES'
$ws.Range("F23").Value = 'Tiene IT'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '8'
$ws.Range("H23").Value = 'ML'
$ws.Range("J23").Value = 'Solo Revisión'

# Row 24
$ws.Range("A24").Value = '0TF26850'
$ws.Range("B24").Value = 'MONTAGNE MASCARILLA BARBIE NIACINAMIDA & COCO'
$ws.Range("C24").Value = 'TRAT.FEMENINO'
$ws.Range("D24").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E24").Value = 'Tiene ES'
$ws.Range("F24").Value = 'Tiene IT'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '1'
$ws.Range("H24").Value = 'UND'
$ws.Range("J24").Value = 'Revisado y Traducido'

# Row 25
$ws.Range("A25").Value = '0TF26851'
$ws.Range("B25").Value = 'MONTAGNE MASCARILLA BARBIE ROSA & VITAMINA E'
$ws.Range("C25").Value = 'TRAT.FEMENINO'
$ws.Range("D25").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E25").Value = 'Tiene ES'
$ws.Range("F25").Value = 'Tiene IT'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '10'
$ws.Range("H25").Value = 'ML'
$ws.Range("J25").Value = 'Revisado y Traducido'

# Row 26
$ws.Range("A26").Value = '0TF26857'
$ws.Range("B26").Value = 'MONTAGNE TUBO PEEL OFF CARBON 50ML'
$ws.Range("C26").Value = 'TRAT.FEMENINO'
$ws.Range("D26").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E26").Value = 'No Tiene ES - TRADUCIDO'
$ws.Range("F26").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '1'
$ws.Range("J26").Value = 'Revisado y Traducido'

# Row 27
$ws.Range("A27").Value = '0TF26858'
$ws.Range("B27").Value = 'MONTAGNE TUBO PEEL OFF ARBOL DE TE 50ML'
$ws.Range("C27").Value = 'TRAT.FEMENINO'
$ws.Range("D27").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E27").Value = 'No Tiene ES - TRADUCIDO'
$ws.Range("F27").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '1'
$ws.Range("J27").Value = 'Revisado y Traducido'

# Row 28
$ws.Range("A28").Value = '0TF26858'
$ws.Range("B28").Value = 'MONTAGNE TUBO PEEL OFF ARBOL DE TE 50ML'
$ws.Range("C28").Value = 'TRAT.FEMENINO'
$ws.Range("D28").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E28").Value = 'Tiene ES'
$ws.Range("F28").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '50'
$ws.Range("H28").Value = 'ML'
$ws.Range("J28").Value = 'Revisado y Traducido'

# Row 29
$ws.Range("A29").Value = '0TF26857'
$ws.Range("B29").Value = 'MONTAGNE TUBO PEEL OFF CARBON 50ML'
$ws.Range("C29").Value = 'TRAT.FEMENINO'
$ws.Range("D29").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E29").Value = 'Tiene ES'
$ws.Range("F29").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '50'
$ws.Range("H29").Value = 'ML'
$ws.Range("J29").Value = 'Revisado y Traducido'

# Row 30
$ws.Range("A30").Value = '0TS04090'
$ws.Range("B30").Value = 'DIOR SOLAR BODY SELF TANING GEL TUBO 150ML'
$ws.Range("C30").Value = 'TRAT.SOLAR'
$ws.Range("D30").Value = 'Tiene PT'
$ws.Range("E30").Value = 'Tiene ES'
$ws.Range("F30").Value = 'Tiene IT'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '150'
$ws.Range("H30").Value = 'ML'
$ws.Range("J30").Value = 'Solo Revisión'

# Row 31
$ws.Range("A31").Value = '0TS04093'
$ws.Range("B31").Value = 'SHISEIDO EXPERT SUN PROTECTOR CLEAN STICK 20G'
$ws.Range("C31").Value = 'TRAT.SOLAR'
$ws.Range("D31").Value = 'Tiene PT'
$ws.Range("E31").Value = 'Tiene ES'
$ws.Range("F31").Value = 'Tiene IT'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '50'
$ws.Range("H31").Value = 'ML'
$ws.Range("J31").Value = 'Solo Revisión'

# Row 32
$ws.Range("A32").Value = '0TS04093'
$ws.Range("B32").Value = 'SHISEIDO EXPERT SUN PROTECTOR CLEAN STICK 20G'
$ws.Range("C32").Value = 'TRAT.SOLAR'
$ws.Range("D32").Value = 'Tiene PT'
$ws.Range("E32").Value = 'Tiene ES'
$ws.Range("F32").Value = 'Tiene IT'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '20'
$ws.Range("H32").Value = 'GR'
$ws.Range("J32").Value = 'Solo Revisión'

# Row 33
$ws.Range("A33").Value = '0TS04092'
$ws.Range("B33").Value = 'SHISEIDO EXPERT SUN PROTECTOR LOTION SPF30 300ML'
$ws.Range("C33").Value = 'TRAT.SOLAR'
$ws.Range("D33").Value = 'Tiene PT'
$ws.Range("E33").Value = 'Tiene ES'
$ws.Range("F33").Value = 'Tiene IT'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '300'
$ws.Range("H33").Value = 'ML'
$ws.Range("J33").Value = 'Solo Revisión'

# Row 34
$ws.Range("A34").Value = '0TF26826'
$ws.Range("B34").Value = 'CLINIQUE MOISTURE SURGE BODY HYDRATOR 200ML'
$ws.Range("C34").Value = 'TRAT.FEMENINO'
$ws.Range("D34").Value = 'Tiene PT'
$ws.Range("E34").Value = 'Tiene ES'
$ws.Range("F34").Value = 'Tiene IT'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '200'
$ws.Range("H34").Value = 'ML'
$ws.Range("J34").Value = 'Solo Revisión'

# Row 35
$ws.Range("A35").Value = '0MO26921'
$ws.Range("B35").Value = 'LANCOME SET MASCARA BIG MOUNSIER'
$ws.Range("C35").Value = 'MAQUILLAJE OJOS'
$ws.Range("D35").Value = 'Tiene PT'
$ws.Range("E35").Value = 'Tiene ES'
$ws.Range("F35").Value = 'Tiene IT'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '1'
$ws.Range("H35").Value = 'UND'
$ws.Range("J35").Value = 'Solo Revisión'

# Row 36
$ws.Range("A36").Value = '0MO26919'
$ws.Range("B36").Value = 'LANCOME MASCARA HYPNOSE DRAMA SET SPRING25'
$ws.Range("C36").Value = 'MAQUILLAJE OJOS'
$ws.Range("D36").Value = 'Tiene PT'
$ws.Range("E36").Value = 'Tiene ES'
$ws.Range("F36").Value = 'Tiene IT'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '1'
$ws.Range("H36").Value = 'UND'
$ws.Range("J36").Value = 'Solo Revisión'

# Row 37
$ws.Range("A37").Value = '0MO27030'
$ws.Range("B37").Value = 'LANCOME SET MASCARA LASH IDOLE 10ML'
$ws.Range("C37").Value = 'MAQUILLAJE OJOS'
$ws.Range("D37").Value = 'Tiene PT'
$ws.Range("E37").Value = 'Tiene ES'
$ws.Range("F37").Value = 'Tiene IT'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '1'
$ws.Range("H37").Value = 'UND'
$ws.Range("J37").Value = 'Solo Revisión'

# Row 38
$ws.Range("A38").Value = '0MO26921'
$ws.Range("B38").Value = 'LANCOME SET MASCARA BIG MOUNSIER'
$ws.Range("C38").Value = 'MAQUILLAJE OJOS'
$ws.Range("D38").Value = 'Tiene PT'
$ws.Range("E38").Value = 'Tiene ES'
$ws.Range("F38").Value = 'Tiene IT'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '1'
$ws.Range("H38").Value = 'UND'
$ws.Range("J38").Value = 'Solo Revisión'

# Row 39
$ws.Range("A39").Value = '6XY00554'
$ws.Range("B39").Value = 'PLUS ONE PRIVATE PLEASURE VIBE'
$ws.Range("C39").Value = 'JUEGOS EROTICOS'
$ws.Range("D39").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E39").Value = 'Tiene ES'
$ws.Range("F39").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '1'
$ws.Range("H39").Value = 'UND'
$ws.Range("J39").Value = 'Revisado y Traducido'

# Row 40
$ws.Range("A40").Value = '6XY00556'
$ws.Range("B40").Value = 'PLUS ONE LUXE MENOPAUSE MASSAGER'
$ws.Range("C40").Value = 'JUEGOS EROTICOS'
$ws.Range("D40").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E40").Value = 'Tiene ES'
$ws.Range("F40").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '1'
$ws.Range("H40").Value = 'UND'
$ws.Range("J40").Value = 'Revisado y Traducido'

# Row 41
$ws.Range("A41").Value = '6XY00555'
$ws.Range("B41").Value = 'PLUS ONE LUXE RIPPLE MULTI VIBE'
$ws.Range("C41").Value = 'JUEGOS EROTICOS'
$ws.Range("D41").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E41").Value = 'Tiene ES'
$ws.Range("F41").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '1'
$ws.Range("H41").Value = 'UND'
$ws.Range("J41").Value = 'Revisado y Traducido'

# Row 42
$ws.Range("A42").Value = '6XS18401'
$ws.Range("B42").Value = 'BLEVIT SF 8 CEREALES Y GALLETA 500GRS + REGALO'
$ws.Range("C42").Value = 'ANEXOS'
$ws.Range("D42").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E42").Value = 'Tiene ES'
$ws.Range("F42").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '1'
$ws.Range("H42").Value = 'UND'
$ws.Range("I42").Value = '"2x Cereales de 250gr cada uno","1x brinde"'
$ws.Range("J42").Value = 'Revisado y Traducido'

# Row 43
$ws.Range("A43").Value = '4EF05448'
$ws.Range("B43").Value = 'KENZO FLOWER SET EDP 100ML + 10ML + CREMA C 75ML'
$ws.Range("C43").Value = 'Set & Pack'
$ws.Range("D43").Value = 'Tiene PT'
$ws.Range("E43").Value = 'Tiene ES'
$ws.Range("F43").Value = 'Tiene IT'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '3'
$ws.Range("H43").Value = 'UND'
$ws.Range("I43").Value = '"1x EDP 100ML","1x Leche Corporal 75ML","1x EDP Mini 10ML"'
$ws.Range("J43").Value = 'Solo Revisión'

# Row 44
$ws.Range("A44").Value = '4EF05447'
$ws.Range("B44").Value = 'KENZO FLOWER SET EDP 50ML + BM75ML + MINI 10ML'
$ws.Range("C44").Value = 'PERF. ESTUCHES MUJER'
$ws.Range("D44").Value = 'Tiene PT'
$ws.Range("E44").Value = 'Tiene ES'
$ws.Range("F44").Value = 'Tiene IT'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '3'
$ws.Range("H44").Value = 'UND'
$ws.Range("J44").Value = 'Solo Revisión'
